$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-14 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-15 Sunday", 2) | Out-Null
$d.Content.Find.Execute("295÷9=32, 7", $true, $false, $false, $false, $false, $true, 1, $false, "313÷4=78, 1", 2) | Out-Null
$d.Content.Find.Execute("354÷9=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "443÷4=110, 3", 2) | Out-Null
$d.Content.Find.Execute("448÷7=64, 0", $true, $false, $false, $false, $false, $true, 1, $false, "249÷2=124, 1", 2) | Out-Null
$d.Content.Find.Execute("765÷7=109, 2", $true, $false, $false, $false, $false, $true, 1, $false, "496÷5=99, 1", 2) | Out-Null
$d.Content.Find.Execute("998÷5=199, 3", $true, $false, $false, $false, $false, $true, 1, $false, "162÷3=54, 0", 2) | Out-Null
$d.Content.Find.Execute("811÷7=115, 6", $true, $false, $false, $false, $false, $true, 1, $false, "881÷7=125, 6", 2) | Out-Null
$d.Content.Find.Execute("950÷6=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "704÷3=234, 2", 2) | Out-Null
$d.Content.Find.Execute("825÷6=137, 3", $true, $false, $false, $false, $false, $true, 1, $false, "804÷4=201, 0", 2) | Out-Null
$d.Content.Find.Execute("921÷5=184, 1", $true, $false, $false, $false, $false, $true, 1, $false, "794÷7=113, 3", 2) | Out-Null
$d.Content.Find.Execute("925÷7=132, 1", $true, $false, $false, $false, $false, $true, 1, $false, "323÷7=46, 1", 2) | Out-Null
$d.Content.Find.Execute("753÷5=150, 3", $true, $false, $false, $false, $false, $true, 1, $false, "861÷6=143, 3", 2) | Out-Null
$d.Content.Find.Execute("739÷8=92, 3", $true, $false, $false, $false, $false, $true, 1, $false, "506÷9=56, 2", 2) | Out-Null
$d.Content.Find.Execute("337÷2=168, 1", $true, $false, $false, $false, $false, $true, 1, $false, "525÷6=87, 3", 2) | Out-Null
$d.Content.Find.Execute("654÷5=130, 4", $true, $false, $false, $false, $false, $true, 1, $false, "249÷8=31, 1", 2) | Out-Null
$d.Content.Find.Execute("236÷5=47, 1", $true, $false, $false, $false, $false, $true, 1, $false, "129÷7=18, 3", 2) | Out-Null
$d.Content.Find.Execute("132÷7=18, 6", $true, $false, $false, $false, $false, $true, 1, $false, "698÷2=349, 0", 2) | Out-Null
$d.Content.Find.Execute("178÷5=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "629÷6=104, 5", 2) | Out-Null
$d.Content.Find.Execute("473÷3=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "955÷5=191, 0", 2) | Out-Null
$d.Content.Find.Execute("400÷7=57, 1", $true, $false, $false, $false, $false, $true, 1, $false, "917÷3=305, 2", 2) | Out-Null
$d.Content.Find.Execute("775÷4=193, 3", $true, $false, $false, $false, $false, $true, 1, $false, "736÷2=368, 0", 2) | Out-Null
$d.Content.Find.Execute("578÷3=192, 2", $true, $false, $false, $false, $false, $true, 1, $false, "198÷8=24, 6", 2) | Out-Null
$d.Content.Find.Execute("323÷5=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "259÷3=86, 1", 2) | Out-Null
$d.Content.Find.Execute("327÷5=65, 2", $true, $false, $false, $false, $false, $true, 1, $false, "386÷7=55, 1", 2) | Out-Null
$d.Content.Find.Execute("533÷7=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "127÷7=18, 1", 2) | Out-Null
$d.Content.Find.Execute("356÷4=89, 0", $true, $false, $false, $false, $false, $true, 1, $false, "992÷6=165, 2", 2) | Out-Null
